$d = $word.ActiveDocument

$d.Content.Find.Execute("Договор аренды нежилого помещения № №23-23с-2ч2", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Договор аренды нежилого помещения № 321", 2)

$d.Content.Find.Execute("г. 2024", $true, $false, $false, $false, $false,
                         $true, 1, $false, "г. 12313", 2)

$d.Content.Find.Execute("«27» ноябрь 2024 г.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "«321» 12 12313 г.", 2)

$d.Content.Find.Execute("Цветков Павел, в лице арендодатель, действующего на основании _____,", $true, $false, $false, $false, $false,
                         $true, 1, $false, "123456 12345, в лице 32, действующего на основании _____,", 2)
